# HeroProfession.xlsx update
# - remove the stale "1=10000" note in D2
# - change the exported field type of atkRatio/armorRatio/maxHPRatio from
#   int32 to number (they are now fractional coefficients, not *10000 ints)
# - rescale the profession-coefficient table (rows 7-15, cols D:F) from
#   "value * 10000" integers down to their real decimal ratios
# - refresh the explanatory cell comment on C3
# - restore the editor's last selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clear the leftover "1=10000" annotation next to the 职业id header.
$ws.Range("D2").ClearContents()

# 2) Field-type row: atkRatio/armorRatio/maxHPRatio switch from int32 to number.
$ws.Range("D6").Value = "number"
$ws.Range("E6").Value = "number"
$ws.Range("F6").Value = "number"

# 3) Rescale the data rows: old values were stored as ratio*10000.
$data = @{
    7  = @{ D = 1.2;  E = 1;    F = 1 }
    8  = @{ D = 0.8;  E = 1.3;  F = 1.3 }
    9  = @{ D = 1.35; E = 0.6;  F = 0.6 }
    10 = @{ D = 1.35; E = 0.6;  F = 0.6 }
    11 = @{ D = 1;    E = 0.6;  F = 0.6 }
    12 = @{ D = 1.25; E = 0.7;  F = 0.7 }
    13 = @{ D = 0.8;  E = 1.25; F = 1.25 }
    14 = @{ D = 1;    E = 1;    F = 1 }
    15 = @{ D = 1.25; E = 0.7;  F = 0.7 }
}
foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
}

# 4) Update the explanatory comment on C3 (keep the bold "Admin:" signature
#    line intact - passing it back through Text() lets the engine re-split
#    it into its own run like the original comment had).
$newComment = "Admin:`n英雄升级属性职业系数表。`n目前此表英雄升级、突破不同养成各属性的职业系数相同。`n比如法师升级和突破的攻击力职业系数相同都是1.35。`n`n只有在att表卡牌升级成长率id（目前id是20）那行填了升级、突破成长率的属性，才会走公式，乘品质、职业系数。`n关联globalconfig表内英雄升级成长率attid"
$ws.Range("C3").Comment.Text($newComment)

# 5) Match the last-saved cell selection.
$ws.Range("L14").Select()
